$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67771723
$ws.Range("C2").Value = 67896547
$ws.Range("D2").Value = 68044296
$ws.Range("E2").Value = 68119911
$ws.Range("F2").Value = 68096423
$ws.Range("G2").Value = 68123067
$ws.Range("H2").Value = 67847796
$ws.Range("I2").Value = 67654279

$ws.Range("B3").Value = 17865082
$ws.Range("C3").Value = 17939648
$ws.Range("D3").Value = 17998603
$ws.Range("E3").Value = 18059917
$ws.Range("F3").Value = 18099311
$ws.Range("G3").Value = 18139809
$ws.Range("H3").Value = 18055158
$ws.Range("I3").Value = 18013939

$ws.Range("B4").Value = 6389901
$ws.Range("C4").Value = 6381639
$ws.Range("D4").Value = 6404542
$ws.Range("E4").Value = 6434640
$ws.Range("F4").Value = 6444785
$ws.Range("G4").Value = 6449822
$ws.Range("H4").Value = 6429482
$ws.Range("I4").Value = 6416656

$ws.Range("B5").Value = 6506418
$ws.Range("C5").Value = 6523269

$ws.Range("D6").Value = 7624798.900000001
$ws.Range("E6").Value = 7634657.890000002
$ws.Range("F6").Value = 7643546.920000001
$ws.Range("G6").Value = 7652586.920000001
$ws.Range("H6").Value = 7661204.900000001
$ws.Range("I6").Value = 7646030.920000001

$ws.Range("B7").Value = 4570466.920000001
$ws.Range("C7").Value = 4602546.880000002
$ws.Range("D7").Value = 4633054.930000001
$ws.Range("E7").Value = 4658936.920000001
$ws.Range("F7").Value = 4682831.920000001
$ws.Range("G7").Value = 4720031.94
$ws.Range("H7").Value = 4739727.920000001
$ws.Range("I7").Value = 4760139.94

$ws.Range("B8").Value = 13694958.92
$ws.Range("C8").Value = 13678443.9
$ws.Range("D8").Value = 13668366.89
$ws.Range("E8").Value = 13636054.93
$ws.Range("F8").Value = 13573442.9
$ws.Range("G8").Value = 13528662.9
$ws.Range("H8").Value = 13426671.92
$ws.Range("I8").Value = 13340651.9

$ws.Range("B9").Value = 5715802.940000001
$ws.Range("C9").Value = 5730861.880000002
$ws.Range("D9").Value = 5750735.930000002
$ws.Range("E9").Value = 5757597.910000001
$ws.Range("F9").Value = 5756688.890000002
$ws.Range("G9").Value = 5765302.930000001
$ws.Range("H9").Value = 5750443.920000001
$ws.Range("I9").Value = 5740710.900000001

$ws.Range("B10").Value = 8988287.92
$ws.Range("C10").Value = 8994098.91
$ws.Range("D10").Value = 8996035.940000001
$ws.Range("E10").Value = 8978336.93

$ws.Range("F11").Value = 7968701.910000001
$ws.Range("G11").Value = 7945532.900000002
$ws.Range("H11").Value = 7888370.900000001
$ws.Range("I11").Value = 7846112.920000002

$ws.Range("B12").Value = 2985001.930000001
$ws.Range("C12").Value = 2974920.930000001
$ws.Range("D12").Value = 2968142.91
$ws.Range("E12").Value = 2959755.9

$ws.Range("F13").Value = 3927102.88
$ws.Range("G13").Value = 3921310.95
$ws.Range("H13").Value = 3896727.920000001
$ws.Range("I13").Value = 3890021.91

$ws.Range("B14").Value = 1055776.95
$ws.Range("C14").Value = 1071096.94
